$d = $word.ActiveDocument

# 1. Reimbursement cap dollar amount: 1000 -> 750
$rng = $d.Content
$rng.Find.Execute("1000", $true, $false, $false, $false, $false, $true, 1, $false, "750", 2)

# 2. Merge "i.e." back into the surrounding sentence (drop the grammar-check
#    markers Word had inserted around it) so the sentence reads as one run:
#    "... items in the same transactions, i.e. soda, gum, candy, etc."
$rng2 = $d.Content
$rng2.Find.Execute(" items in the same transactions, i.e. soda, gum, candy, etc.", $true, $false, $false, $false, $false, $true, 1, $false, " items in the same transactions, i.e. soda, gum, candy, etc.", 2)

# 3. Re-key the phone number text so the stale page-break bookkeeping before
#    "8803" is recalculated away.
$rng3 = $d.Content
$rng3.Find.Execute("8803", $true, $false, $false, $false, $false, $true, 1, $false, "8803", 2)
